# Commit: "excel file added 123"
# Reproduce the author's edit: type the number 123 into cell A1 of Sheet1,
# press Enter (which leaves the active selection on A2), and set the sheet's
# print setup to A4/portrait - matching the saved worksheet state in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the value into A1
$ws.Range("A1").Value = 123

# Page setup, as reflected in the saved file (paperSize=9 -> A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# After typing into A1 and pressing Enter, the active cell moves to A2
$ws.Range("A2").Select() | Out-Null
